# Data update: increment a handful of enrollment/payment counts in the
# "Resumo Inscrições Subsequente" sheet, matching the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Agropecuária - Bambuí): Inscritos, Pagos, Inscrições homologadas
$ws.Range("E2").Value = 31
$ws.Range("F2").Value = 17
$ws.Range("H2").Value = 19

# Row 10 (Eletrotécnica - Conselheiro Lafaiete): Inscritos
$ws.Range("E10").Value = 36

# Row 14 (Técnico Subsequente - Conselheiro Lafaiete group): Pagos, Inscrições homologadas
$ws.Range("F14").Value = 19
$ws.Range("H14").Value = 21

# Row 15: Inscritos
$ws.Range("E15").Value = 110
